$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Seed the new "dne" shared string first so it lands at the same
# sharedStrings index the target workbook uses, before the other
# brand-new key strings get appended.
$ws.Range("C3").Value = "dne"

$ws.Range("A2").Value = "GREENKEY"
$ws.Range("B2").Value = 100
$ws.Range("C2").Value = 85

$ws.Range("A3").Value = "PURPLEKEY"
$ws.Range("B3").Value = 98.1

$ws.Range("A4").Value = "REDKEY"
$ws.Range("B4").Value = 100
$ws.Range("C4").Value = 100

$ws.Range("A5").Value = "ALPHAKEY"
$ws.Range("B5").Value = 98.55
$ws.Range("C5").Value = "dne"

$ws.Range("A6").Value = "ACETAKEY"
$ws.Range("B6").Value = 100
$ws.Range("C6").Value = "dne"

$ws.Range("A7").Value = "BETAKEY"
$ws.Range("B7").Value = 98.6
$ws.Range("C7").Value = 54

$ws.Range("A8").Value = "CHARLIEKEY"
$ws.Range("B8").Value = 88
$ws.Range("C8").Value = 79

[void]$ws.Range("C9").Select()
